$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 07:49 AM"

# --- "1 Month Performance" sheet: refresh stock list/values ---
$perf = $wb.Worksheets.Item("1 Month Performance")
$perf.Range("C8").Value = 65.36799999999999
$perf.Range("C10").Value = 52.2777
$perf.Range("C12").Value = 44.3722
$perf.Range("C13").Value = 40.177
$perf.Range("B15").Value = "SEJALLTD"
$perf.Range("C15").Value = 39.4454
$perf.Range("B16").Value = "SANDUMA"
$perf.Range("C16").Value = 38.8163
$perf.Range("C17").Value = 37.8136
$perf.Range("B18").Value = "SOUTHBANK"
$perf.Range("C18").Value = 37.1844
$perf.Range("B19").Value = "RAMAPHO"
$perf.Range("C19").Value = 36.9374
$perf.Range("C20").Value = 36.713
$perf.Range("C21").Value = 36.1733
$perf.Range("C22").Value = 35.7172
$perf.Range("C24").Value = 35.2478
$perf.Range("C25").Value = 33.9051
$perf.Range("C26").Value = 33.7654
$perf.Range("C28").Value = 32.5448
$perf.Range("B29").Value = "TATVA"
$perf.Range("C29").Value = 31.5056
$perf.Range("B30").Value = "ATHERENERG"
$perf.Range("C30").Value = 31.182
$perf.Range("B31").Value = "TARACHAND"
$perf.Range("C31").Value = 31.0974
$perf.Range("C34").Value = 27.7454
$perf.Range("C36").Value = 27.3706
$perf.Range("C37").Value = 27.3461
$perf.Range("C39").Value = 26.8616
$perf.Range("C41").Value = 26.0381
$perf.Range("C42").Value = 25.8512
$perf.Range("C43").Value = 25.4212
$perf.Range("C44").Value = 24.8065
$perf.Range("C45").Value = 24.1422
$perf.Range("C47").Value = 23.9933
$perf.Range("B48").Value = "SKYGOLD"
$perf.Range("C48").Value = 23.9127
$perf.Range("B49").Value = "AUBANK"
$perf.Range("C49").Value = 23.6754
$perf.Range("B50").Value = "DCBBANK"
$perf.Range("C50").Value = 23.6742
$perf.Range("B51").Value = "PRECWIRE"
$perf.Range("C51").Value = 23.1358
$perf.Range("B52").Value = "TDPOWERSYS"
$perf.Range("C52").Value = 22.9845
$perf.Range("B53").Value = "ETHOSLTD"
$perf.Range("C53").Value = 22.9743
$perf.Range("B55").Value = "GUJTHEM"
$perf.Range("C55").Value = 22.3332
$perf.Range("B56").Value = "INDIANB"
$perf.Range("C56").Value = 22.2382
$perf.Range("B57").Value = "ORBTEXP"
$perf.Range("C57").Value = 21.6706
$perf.Range("C58").Value = 21.0713
$perf.Range("C59").Value = 20.6412
$perf.Range("C60").Value = 20.1086
$perf.Range("B61").Value = "ATL"
$perf.Range("C61").Value = 19.8652
$perf.Range("B62").Value = "RBLBANK"
$perf.Range("C62").Value = 19.7715
$perf.Range("B63").Value = "GRMOVER"
$perf.Range("C63").Value = 19.5401
$perf.Range("B64").Value = "FEDERALBNK"
$perf.Range("C64").Value = 19.5343
$perf.Range("C65").Value = 19.4932
$perf.Range("B66").Value = "CEATLTD"
$perf.Range("C66").Value = 19.4529
$perf.Range("C72").Value = 18.9612
$perf.Range("C73").Value = 18.5125
$perf.Range("C74").Value = 18.2884

"Update complete"
